# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the
# "6489d606-5d3b-409c-bc93-e42c01cf99eb" file now that a new xliff
# round-trip has completed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 6489d606 row.
$wsOverview.Range("G3").Value = "2016-09-01 06:55:25"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 6489d606 row.
$wsZhCn.Range("H3").Value = "2016-09-01 06:55:20"
$wsZhCn.Range("K3").Value = "2016-09-01 06:55:37"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 6489d606 row.
$wsDeDe.Range("H3").Value = "2016-09-01 06:55:25"
$wsDeDe.Range("K3").Value = "2016-09-01 06:55:44"
